# Final compile tweaks for Fall 2020 class:
#  - Fix the misspelled sheet name "Enrolment Statistics" -> "Enrollment Statistics"
#  - Leave that sheet active/selected as the last-used tab before saving

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Enrolment Statistics")
$ws.Name = "Enrollment Statistics"

# Make the renamed sheet the active tab (matches tabSelected/activeTab moving
# to this sheet in the saved workbook).
$ws.Activate()
